$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E4").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("B6").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("E6").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("F6").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("B10").Value = "Rogério-Elem. Máquinas"
$ws.Range("D10").Value = "[Paulo Rob.-Usin. CNC-3A, Wellington-Trat. Termicos-3A, Paulo Rob.-M.A.Comp.CAD / CAM-3A, Joel L.-Tec. Fundição-3A]"
$ws.Range("E10").Value = "-"
$ws.Range("F10").Value = "[Tiago Freitas-M.S.R. ar Cond.-3A, Gisele-Ens. Dest. Não Desti.-3A, Ivan-Tec. Soldagem-3A, Aselmo-M. Motor Endot.-3A]"
$ws.Range("B11").Value = "[Humberto-C.pneumática-3A, Ludoff-Comandos Eletricos-3A, Valmir-Calderaria-3A, Ludoff-Cont.Lóg.Prog CLP-3A]"
$ws.Range("C11").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("D11").Value = "[Paulo Rob.-Usin. CNC-3A, Wellington-Trat. Termicos-3A, Paulo Rob.-M.A.Comp.CAD / CAM-3A, Joel L.-Tec. Fundição-3A]"
$ws.Range("E11").Value = "-"
$ws.Range("F11").Value = "[Leonardo-Retífica-3A, Aderci-Fresagem-3A, Nilton-Metrologia 2-3A, Leonardo-Mec. Manut.Equip. ind.-3A]"
$ws.Range("B12").Value = "[Humberto-C.pneumática-3A, Ludoff-Comandos Eletricos-3A, Valmir-Calderaria-3A, Ludoff-Cont.Lóg.Prog CLP-3A]"
$ws.Range("C12").Value = "-"
$ws.Range("D12").Value = "[Paulo Rob.-Usin. CNC-3A, Wellington-Trat. Termicos-3A, Paulo Rob.-M.A.Comp.CAD / CAM-3A, Joel L.-Tec. Fundição-3A]"
$ws.Range("E12").Value = "-"
$ws.Range("F12").Value = "[Leonardo-Retífica-3A, Aderci-Fresagem-3A, Nilton-Metrologia 2-3A, Leonardo-Mec. Manut.Equip. ind.-3A]"
$ws.Range("B14").Value = "[Humberto-C.pneumática-3A, Ludoff-Comandos Eletricos-3A, Valmir-Calderaria-3A, Ludoff-Cont.Lóg.Prog CLP-3A]"
$ws.Range("D14").Value = "[Tiago Freitas-M.S.R. ar Cond.-3A, Gisele-Ens. Dest. Não Desti.-3A, Ivan-Tec. Soldagem-3A, Aselmo-M. Motor Endot.-3A]"
$ws.Range("E14").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("F14").Value = "[Leonardo-Retífica-3A, Aderci-Fresagem-3A, Nilton-Metrologia 2-3A, Leonardo-Mec. Manut.Equip. ind.-3A]"
$ws.Range("B15").Value = "[Humberto-C.pneumática-3A, Ludoff-Comandos Eletricos-3A, Valmir-Calderaria-3A, Ludoff-Cont.Lóg.Prog CLP-3A]"
$ws.Range("D15").Value = "[Tiago Freitas-M.S.R. ar Cond.-3A, Gisele-Ens. Dest. Não Desti.-3A, Ivan-Tec. Soldagem-3A, Aselmo-M. Motor Endot.-3A]"
$ws.Range("E15").Value = "-"
$ws.Range("F15").Value = "[Leonardo-Retífica-3A, Aderci-Fresagem-3A, Nilton-Metrologia 2-3A, Leonardo-Mec. Manut.Equip. ind.-3A]"
$ws.Range("B16").Value = "Rogério-Elem. Máquinas"
$ws.Range("D16").Value = "[Joel L.-Tec. Fundição-3A, Paulo Rob.-M.A.Comp.CAD / CAM-3A, Wellington-Trat. Termicos-3A, Paulo Rob.-Usin. CNC-3A]"
$ws.Range("E16").Value = "-"
$ws.Range("F16").Value = "[Ivan-Tec. Soldagem-3A, Aselmo-M. Motor Endot.-3A, Gisele-Ens. Dest. Não Desti.-3A, Tiago Freitas-M.S.R. ar Cond.-3A]"
$ws.Range("E18").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("C19").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
$ws.Range("E21").Value = "[Gisele-Ens. Dest. Não Desti.-3A, -, -, -]"
